$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("משתמשים")
$ws.Activate()

# Remove the three broken/incomplete registration rows (old rows 7, 8, 9 -
# each a leftover "Yoni" test entry) so the remaining real records shift up.
$ws.Rows.Item(7).Resize(3).EntireRow.Delete() | Out-Null

# The row that used to be row 12 ("Sap" placeholder row with a count of 1)
# is now free at the bottom (row 12) - fill it in with the finished
# registration test record.
$ws.Range("A12").Value = "test"
$ws.Range("B12").Value = "test"
$ws.Range("C12").Value = "test"
$ws.Range("D12").Value = "test"
$ws.Range("F12").Value = "test"
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = $false

$ws.Range("M10").Select()
